$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.110744889959449
$ws.Range("C2").Value = 0.2370466469056112
$ws.Range("D2").Value = 0.07896879185899763
$ws.Range("E2").Value = 0.1041169720498214
$ws.Range("G2").Value = 0.5314283454143123
$ws.Range("H2").Value = 0.649162030729812
$ws.Range("L2").Value = 0.2078881262378616
$ws.Range("N2").Value = 1.109682608641549
$ws.Range("O2").Value = 2.330101827823086

$ws.Range("B3").Value = 1.007981929872017
$ws.Range("C3").Value = 0.2281235549205718
$ws.Range("D3").Value = 0.07157096782734129
$ws.Range("E3").Value = 0.1047645686840024
$ws.Range("G3").Value = 0.5260900438135394
$ws.Range("H3").Value = 0.6510438205903029
$ws.Range("L3").Value = 0.1985677997574413
$ws.Range("N3").Value = 1.118982550764215
$ws.Range("O3").Value = 2.322401474574406

$ws.Range("B4").Value = 0.9450319952749737
$ws.Range("C4").Value = 0.2226116043022586
$ws.Range("D4").Value = 0.06706314025053928
$ws.Range("E4").Value = 0.1052180906332225
$ws.Range("G4").Value = 0.5232431947960379
$ws.Range("H4").Value = 0.652553346698312
$ws.Range("L4").Value = 0.192948106710503
$ws.Range("N4").Value = 1.125157048392737
$ws.Range("O4").Value = 2.319316495077544

$ws.Range("B5").Value = 0.9194177838570567
$ws.Range("C5").Value = 0.2203572743890874
$ws.Range("D5").Value = 0.06523482299692773
$ws.Range("E5").Value = 0.1054169606215254
$ws.Range("G5").Value = 0.5221911962112529
$ws.Range("H5").Value = 0.6532575216318008
$ws.Range("L5").Value = 0.1906840168717281
$ws.Range("N5").Value = 1.127790035911445
$ws.Range("O5").Value = 2.318471975254766

$ws.Range("B6").Value = 0.9151669254649732
$ws.Range("C6").Value = 0.2199824563851394
$ws.Range("D6").Value = 0.06493175507553417
$ws.Range("E6").Value = 0.1054508318127407
$ws.Range("G6").Value = 0.5220230356730298
$ws.Range("H6").Value = 0.6533798262938717
$ws.Range("L6").Value = 0.1903096378466103
$ws.Range("N6").Value = 1.128234300106968
$ws.Range("O6").Value = 2.318356650041295

$ws.Range("B7").Value = 0.9446863956681
$ws.Range("C7").Value = 0.2225812344594402
$ws.Range("D7").Value = 0.06703844787510604
$ws.Range("E7").Value = 0.105220715749514
$ws.Range("G7").Value = 0.5232285697327654
$ws.Range("H7").Value = 0.6525624829803149
$ws.Range("L7").Value = 0.1929174671089413
$ws.Range("N7").Value = 1.125192084638783
$ws.Range("O7").Value = 2.319303435552484

$ws.Range("B8").Value = 1.075282560328162
$ws.Range("C8").Value = 0.2339769396381399
$ws.Range("D8").Value = 0.07641086174832878
$ws.Range("E8").Value = 0.104328662777359
$ws.Range("G8").Value = 0.5294981231249665
$ws.Range("H8").Value = 0.649737370894087
$ws.Range("L8").Value = 0.2046531359562778
$ws.Range("N8").Value = 1.112792910198969
$ws.Range("O8").Value = 2.327105399492808

$ws.Range("B9").Value = 1.3324964241462
$ws.Range("C9").Value = 0.2560546382573818
$ws.Range("D9").Value = 0.0950649385653719
$ws.Range("E9").Value = 0.1030228709396273
$ws.Range("G9").Value = 0.5452237253171859
$ws.Range("H9").Value = 0.6470080739191388
$ws.Range("L9").Value = 0.2284823412919792
$ws.Range("N9").Value = 1.092158909350879
$ws.Range("O9").Value = 2.355470420149203

$ws.Range("B10").Value = 1.522102256486676
$ws.Range("C10").Value = 0.2721039768154867
$ws.Range("D10").Value = 0.1089408607896161
$ws.Range("E10").Value = 0.1023340002117727
$ws.Range("G10").Value = 0.5588878932375394
$ws.Range("H10").Value = 0.6467188476803187
$ws.Range("L10").Value = 0.2464864276758902
$ws.Range("N10").Value = 1.079238967191756
$ws.Range("O10").Value = 2.384321043778669

$ws.Range("B11").Value = 1.608486632717643
$ws.Range("C11").Value = 0.2793666834123201
$ws.Range("D11").Value = 0.1152912297184798
$ws.Range("E11").Value = 0.1020793815942174
$ws.Range("G11").Value = 0.5655665551803821
$ws.Range("H11").Value = 0.6469604884640887
$ws.Range("L11").Value = 0.2547848567492537
$ws.Range("N11").Value = 1.073846980146946
$ws.Range("O11").Value = 2.399195675578341

$ws.Range("B12").Value = 1.641215887461442
$ws.Range("C12").Value = 0.2821112360297775
$ws.Range("D12").Value = 0.1177014626576067
$ws.Range("E12").Value = 0.1019914142280633
$ws.Range("G12").Value = 0.5681624274364481
$ws.Range("H12").Value = 0.6471056935015582
$ws.Range("L12").Value = 0.2579427823153111
$ws.Range("N12").Value = 1.071874924022204
$ws.Range("O12").Value = 2.405080707794014

$ws.Range("B13").Value = 1.634166309548107
$ws.Range("C13").Value = 0.281520402867784
$ws.Range("D13").Value = 0.1171821318632027
$ws.Range("E13").Value = 0.1020099836632191
$ws.Range("G13").Value = 0.5676003844228035
$ws.Range("H13").Value = 0.647072032098805
$ws.Range("L13").Value = 0.2572619788692379
$ws.Range("N13").Value = 1.072296538743913
$ws.Range("O13").Value = 2.403802028182156

$ws.Range("B14").Value = 1.611178949176917
$ws.Range("C14").Value = 0.279592594213284
$ws.Range("D14").Value = 0.1154894110796363
$ws.Range("E14").Value = 0.1020719750764627
$ws.Range("G14").Value = 0.5657787787200874
$ws.Range("H14").Value = 0.6469713581307133
$ws.Range("L14").Value = 0.2550443511906479
$ws.Range("N14").Value = 1.073683339668328
$ws.Range("O14").Value = 2.399674779980643

$ws.Range("B15").Value = 1.597100742066573
$ws.Range("C15").Value = 0.2784110111072096
$ws.Range("D15").Value = 0.1144532851933491
$ws.Range("E15").Value = 0.1021110472797648
$ws.Range("G15").Value = 0.5646717007809343
$ws.Range("H15").Value = 0.6469166868619283
$ws.Range("L15").Value = 0.2536880051599866
$ws.Range("N15").Value = 1.074541881134031
$ws.Range("O15").Value = 2.397179601176134

$ws.Range("B16").Value = 1.516459350270907
$ws.Range("C16").Value = 0.2716285573216908
$ws.Range("D16").Value = 0.1085266168871044
$ws.Range("E16").Value = 0.1023518226108475
$ws.Range("G16").Value = 0.558460761475132
$ws.Range("H16").Value = 0.646710568719044
$ws.Range("L16").Value = 0.2459462778580246
$ws.Range("N16").Value = 1.079601110928493
$ws.Range("O16").Value = 2.38338422726406

$ws.Range("B17").Value = 1.467021085190595
$ws.Range("C17").Value = 0.2674578240242056
$ws.Range("D17").Value = 0.1049005680954309
$ws.Range("E17").Value = 0.1025145802274832
$ws.Range("G17").Value = 0.554769263794725
$ws.Range("H17").Value = 0.6466797341213777
$ws.Range("L17").Value = 0.2412246551525499
$ws.Range("N17").Value = 1.082829079718543
$ws.Range("O17").Value = 2.375369968728933

$ws.Range("B18").Value = 1.438597995977716
$ws.Range("C18").Value = 0.2650553420652102
$ws.Range("D18").Value = 0.1028185508424428
$ws.Range("E18").Value = 0.102613723972711
$ws.Range("G18").Value = 0.5526895537950907
$ws.Range("H18").Value = 0.6466971252599762
$ws.Range("L18").Value = 0.238519099260003
$ws.Range("N18").Value = 1.084731407299984
$ws.Range("O18").Value = 2.370925091465324

$ws.Range("B19").Value = 1.428976620178275
$ws.Range("C19").Value = 0.2642412924998325
$ws.Range("D19").Value = 0.1021142321303188
$ws.Range("E19").Value = 0.1026482420193151
$ws.Range("G19").Value = 0.5519928700392427
$ws.Range("H19").Value = 0.6467090455726066
$ws.Range("L19").Value = 0.2376047982322973
$ws.Range("N19").Value = 1.085383349974073
$ws.Range("O19").Value = 2.369448402497795

$ws.Range("B20").Value = 1.472282591331577
$ws.Range("C20").Value = 0.2679021782061284
$ws.Range("D20").Value = 0.1052861956072633
$ws.Range("E20").Value = 0.1024966820612718
$ws.Range("G20").Value = 0.5551577216983929
$ws.Range("H20").Value = 0.6466793808233575
$ws.Range("L20").Value = 0.2417262252192813
$ws.Range("N20").Value = 1.082480728842732
$ws.Range("O20").Value = 2.37620604851395

$ws.Range("B21").Value = 1.6179304353725
$ws.Range("C21").Value = 0.2801589938524955
$ws.Range("D21").Value = 0.1159864553102068
$ws.Range("E21").Value = 0.1020535373118214
$ws.Range("G21").Value = 0.5663120132326895
$ws.Range("H21").Value = 0.6469994708193383
$ws.Range("L21").Value = 0.2556953024242716
$ws.Range("N21").Value = 1.073274109168025
$ws.Range("O21").Value = 2.400880199553001

$ws.Range("B22").Value = 1.713220215657373
$ws.Range("C22").Value = 0.2881363594351001
$ws.Range("D22").Value = 0.1230116705564086
$ws.Range("E22").Value = 0.1018131777778883
$ws.Range("G22").Value = 0.5739915165035114
$ws.Range("H22").Value = 0.6475216938885495
$ws.Range("L22").Value = 0.2649151556299785
$ws.Range("N22").Value = 1.06766373095914
$ws.Range("O22").Value = 2.418477260742492

$ws.Range("B23").Value = 1.662353591349074
$ws.Range("C23").Value = 0.28388178417066
$ws.Range("D23").Value = 0.1192592548378144
$ws.Range("E23").Value = 0.1019369538961676
$ws.Range("G23").Value = 0.5698570925394364
$ws.Range("H23").Value = 0.647214320589427
$ws.Range("L23").Value = 0.2599861146824907
$ws.Range("N23").Value = 1.070620890088968
$ws.Range("O23").Value = 2.408950564999259

$ws.Range("B24").Value = 1.469903865444223
$ws.Range("C24").Value = 0.2677013002190165
$ws.Range("D24").Value = 0.1051118451675563
$ws.Range("E24").Value = 0.1025047564656418
$ws.Range("G24").Value = 0.5549819672623642
$ws.Range("H24").Value = 0.6466794311577218
$ws.Range("L24").Value = 0.2414994374541379
$ws.Range("N24").Value = 1.082638073345713
$ws.Range("O24").Value = 2.37582755031454

$ws.Range("B25").Value = 1.262798954995787
$ws.Range("C25").Value = 0.2501115873144641
$ws.Range("D25").Value = 0.08998875932674366
$ws.Range("E25").Value = 0.1033286286602682
$ws.Range("G25").Value = 0.5406002586178573
$ws.Range("H25").Value = 0.6474452561602959
$ws.Range("L25").Value = 0.2219486400094866
$ws.Range("N25").Value = 1.097347339628051
$ws.Range("O25").Value = 2.346393447255934
